# Refresh Pais worksheet: updated country case counts (and a couple of
# re-sorted country-name cells) plus the "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 23 de Marzo de 2020 a las 01:16'

$ws.Range("B6").Value = 32783
$ws.Range("C6").Value = 8576
$ws.Range("E6").Value = 32189
$ws.Range("G6").Value = 114
$ws.Range("H6").Value = 416

$ws.Range("B21").Value = 1470
$ws.Range("C21").Value = 142
$ws.Range("E21").Value = 1436

$ws.Range("A54").Value = 'Panama'
$ws.Range("B54").Value = 313
$ws.Range("C54").Value = 68
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 309
$ws.Range("F54").Value = 7
$ws.Range("H54").Value = 3

$ws.Range("A55").Value = 'Sudafrica'
$ws.Range("B55").Value = 274
$ws.Range("C55").Value = 34
$ws.Range("D55").Value = 2
$ws.Range("E55").Value = 272
$ws.Range("F55").Value = 0
$ws.Range("H55").Value = 0

$ws.Range("A56").Value = 'Argentina'
$ws.Range("B56").Value = 266
$ws.Range("C56").Value = 108
$ws.Range("D56").Value = 27
$ws.Range("E56").Value = 235
$ws.Range("F56").Value = 0
$ws.Range("H56").Value = 4

$ws.Range("A57").Value = 'Croacia'
$ws.Range("B57").Value = 254
$ws.Range("C57").Value = 48
$ws.Range("D57").Value = 5
$ws.Range("E57").Value = 248
$ws.Range("F57").Value = 5
$ws.Range("H57").Value = 1

$ws.Range("A58").Value = 'Mexico'
$ws.Range("B58").Value = 251
$ws.Range("C58").Value = 48
$ws.Range("D58").Value = 4
$ws.Range("E58").Value = 245
$ws.Range("F58").Value = 1
$ws.Range("H58").Value = 2

$ws.Range("A59").Value = 'Libano'
$ws.Range("B59").Value = 248
$ws.Range("C59").Value = 18
$ws.Range("D59").Value = 8
$ws.Range("E59").Value = 236
$ws.Range("F59").Value = 4
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 4

$ws.Range("A60").Value = 'Irak'
$ws.Range("B60").Value = 233
$ws.Range("C60").Value = 19
$ws.Range("D60").Value = 57
$ws.Range("E60").Value = 156
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 20

$ws.Range("A61").Value = 'Colombia'
$ws.Range("B61").Value = 231
$ws.Range("C61").Value = 35
$ws.Range("D61").Value = 3
$ws.Range("E61").Value = 226
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 2

$ws.Range("A71").Value = 'Uruguay'
$ws.Range("B71").Value = 158
$ws.Range("C71").Value = 48
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 158
$ws.Range("H71").Value = 0

$ws.Range("A72").Value = 'Emiratos Arabes Unidos'
$ws.Range("B72").Value = 153
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 38
$ws.Range("E72").Value = 113
$ws.Range("F72").Value = 2
$ws.Range("H72").Value = 2

$ws.Range("A73").Value = 'Lituania'
$ws.Range("B73").Value = 143
$ws.Range("C73").Value = 44
$ws.Range("D73").Value = 1
$ws.Range("E73").Value = 141
$ws.Range("F73").Value = 1
$ws.Range("H73").Value = 1

$ws.Range("A74").Value = 'Letonia'
$ws.Range("B74").Value = 139
$ws.Range("C74").Value = 15
$ws.Range("D74").Value = 1
$ws.Range("E74").Value = 138
$ws.Range("F74").Value = 0
$ws.Range("H74").Value = 0

$ws.Range("A75").Value = 'Costa Rica'
$ws.Range("B75").Value = 134
$ws.Range("C75").Value = 17
$ws.Range("D75").Value = 2
$ws.Range("E75").Value = 130
$ws.Range("F75").Value = 2
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 2

$ws.Range("A76").Value = 'Hungria'
$ws.Range("B76").Value = 131
$ws.Range("C76").Value = 28
$ws.Range("D76").Value = 16
$ws.Range("E76").Value = 109
$ws.Range("F76").Value = 6
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 6

$ws.Range("A84").Value = 'Nueva Zelanda'
$ws.Range("B84").Value = 102
$ws.Range("C84").Value = 50
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 102
$ws.Range("F84").Value = 0
$ws.Range("H84").Value = 0

$ws.Range("A85").Value = 'Republica de Chipre'
$ws.Range("B85").Value = 95
$ws.Range("C85").Value = 11
$ws.Range("D85").Value = 3
$ws.Range("E85").Value = 91

$ws.Range("A86").Value = 'Moldavia'
$ws.Range("B86").Value = 94
$ws.Range("C86").Value = 14
$ws.Range("D86").Value = 1
$ws.Range("E86").Value = 92
$ws.Range("F86").Value = 3
$ws.Range("H86").Value = 1

$ws.Range("A87").Value = 'Malta'
$ws.Range("B87").Value = 90
$ws.Range("C87").Value = 17
$ws.Range("E87").Value = 88
$ws.Range("F87").Value = 1
$ws.Range("H87").Value = 0

$ws.Range("A88").Value = 'Albania'
$ws.Range("B88").Value = 89
$ws.Range("C88").Value = 13
$ws.Range("E88").Value = 85
$ws.Range("H88").Value = 2

$ws.Range("A89").Value = 'Brunei'
$ws.Range("B89").Value = 88
$ws.Range("C89").Value = 5
$ws.Range("E89").Value = 86
$ws.Range("F89").Value = 2

$ws.Range("A90").Value = 'Camboya'
$ws.Range("B90").Value = 84
$ws.Range("C90").Value = 31
$ws.Range("D90").Value = 2
$ws.Range("E90").Value = 82
$ws.Range("F90").Value = 0

$ws.Range("A91").Value = 'Sri Lanka'
$ws.Range("B91").Value = 82
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 3
$ws.Range("E91").Value = 79
$ws.Range("F91").Value = 2

$ws.Range("A92").Value = 'Bielorrusia'
$ws.Range("B92").Value = 76
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 15
$ws.Range("E92").Value = 61
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0

$ws.Range("A93").Value = 'Tunez'
$ws.Range("C93").Value = 15
$ws.Range("D93").Value = 1
$ws.Range("E93").Value = 71
$ws.Range("F93").Value = 7
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 3

$ws.Range("A94").Value = 'Burkina Faso'
$ws.Range("B94").Value = 75
$ws.Range("C94").Value = 11
$ws.Range("D94").Value = 5
$ws.Range("E94").Value = 66
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 4

$ws.Range("A95").Value = 'Ucrania'
$ws.Range("B95").Value = 73
$ws.Range("C95").Value = 26
$ws.Range("D95").Value = 1
$ws.Range("E95").Value = 69
$ws.Range("F95").Value = 0
$ws.Range("H95").Value = 3

$ws.Range("A96").Value = 'Venezuela'
$ws.Range("B96").Value = 70
$ws.Range("C96").Value = 0
$ws.Range("D96").Value = 15
$ws.Range("E96").Value = 55
$ws.Range("F96").Value = 2

$ws.Range("A97").Value = 'Senegal'
$ws.Range("B97").Value = 67
$ws.Range("C97").Value = 11
$ws.Range("D97").Value = 5
$ws.Range("E97").Value = 62

$ws.Range("A119").Value = 'Ghana'

$ws.Range("A120").Value = 'Puerto Rico'

$ws.Range("A128").Value = 'Polinesia Francesa'
$ws.Range("C128").Value = 3

$ws.Range("A129").Value = 'Guayana Francesa'
$ws.Range("E129").Value = 18
$ws.Range("H129").Value = 0

$ws.Range("A130").Value = 'Guyana'
$ws.Range("B130").Value = 18
$ws.Range("E130").Value = 17
$ws.Range("H130").Value = 1

$ws.Range("A131").Value = 'Togo'
$ws.Range("B131").Value = 16
$ws.Range("E131").Value = 16

$ws.Range("A139").Value = 'Etiopia'
$ws.Range("C139").Value = 2

$ws.Range("A140").Value = 'Mayotte'
$ws.Range("C140").Value = 0

$ws.Range("B142").Value = 9
$ws.Range("C142").Value = 4
$ws.Range("E142").Value = 8

$ws.Range("A144").Value = 'Islas Virgenes de los Estados Unidos'

$ws.Range("A145").Value = 'Guinea Ecuatorial'

$ws.Range("A146").Value = 'Bermudas'
$ws.Range("B146").Value = 6
$ws.Range("C146").Value = 4
$ws.Range("E146").Value = 6

$ws.Range("A149").Value = 'Surinam'
$ws.Range("E149").Value = 5
$ws.Range("H149").Value = 0

$ws.Range("A150").Value = 'Gabon'
$ws.Range("B150").Value = 5
$ws.Range("H150").Value = 1

$ws.Range("A154").Value = 'Nueva Caledonia'
$ws.Range("B154").Value = 4
$ws.Range("E154").Value = 4

$ws.Range("A155").Value = 'San Bartolome'

$ws.Range("A156").Value = 'Republica de Africa Central'
$ws.Range("C156").Value = 0

$ws.Range("A157").Value = 'Zambia'
$ws.Range("C157").Value = 1

$ws.Range("A158").Value = 'Liberia'

$ws.Range("A159").Value = 'Congo'

$ws.Range("A160").Value = 'El Salvador'

$ws.Range("A161").Value = 'Namibia'

$ws.Range("A162").Value = 'Cabo Verde'

$ws.Range("A163").Value = 'Madagascar'

$ws.Range("A164").Value = 'Zimbabue'
$ws.Range("E164").Value = 3
$ws.Range("H164").Value = 0

$ws.Range("A165").Value = 'Curazao'

$ws.Range("A166").Value = 'Islas Caimanes'
$ws.Range("B166").Value = 3
$ws.Range("H166").Value = 1

$ws.Range("A167").Value = 'Mauritania'

$ws.Range("A168").Value = 'Benin'

$ws.Range("A169").Value = 'Guinea'
$ws.Range("C169").Value = 0

$ws.Range("A170").Value = 'Niger'
$ws.Range("C170").Value = 1

$ws.Range("A171").Value = 'Butan'

$ws.Range("A172").Value = 'Nicaragua'

$ws.Range("A173").Value = 'Santa Lucia'

$ws.Range("A174").Value = 'Haiti'

$ws.Range("A175").Value = 'Angola'

$ws.Range("A180").Value = 'Montserrat'

$ws.Range("A181").Value = 'Dominica'
$ws.Range("C181").Value = 1

$ws.Range("A182").Value = 'Republica de Yibuti'
$ws.Range("C182").Value = 0

$ws.Range("A183").Value = 'Republica del Chad'

$ws.Range("A184").Value = 'Santa Sede'

$ws.Range("A185").Value = 'San Vicente y las Granadinas'

$ws.Range("A186").Value = 'Siria'
$ws.Range("C186").Value = 1

$ws.Range("A187").Value = 'Granada'

$ws.Range("A188").Value = 'Papua Nueva Guinea'
$ws.Range("C188").Value = 0

$ws.Range("A189").Value = 'Antigua y Barbuda'

$ws.Range("A190").Value = 'Timor Oriental'

$ws.Range("A191").Value = 'Eritrea'

$ws.Range("A192").Value = 'Uganda'

$ws.Range("A193").Value = 'Gambia'

$ws.Range("A194").Value = 'Somalia'
